# Word doc edit: split the "MVN project (...)" TODO item into two bullets.
# 1) The existing paragraph's text is replaced by a new note about switching
#    to the Plot tab after drawing a plot (the hidden _GoBack bookmark that
#    sits between its two runs is left untouched).
# 2) A brand-new paragraph is appended right after it, re-stating the
#    original "MVN project (...)" text split across two runs.

$d = $word.ActiveDocument

# Non-breaking space used in the original Slovak text ("a<nbsp>spol.")
$nbsp = [char]0xA0

# --- Locate the paragraph that currently holds the "MVN project" text ---
$targetIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.StartsWith("MVN project (")) {
        $targetIdx = $i
    }
}

if ($targetIdx -eq -1) {
    throw "Could not find the 'MVN project (...)' paragraph"
}

# --- Step 1: rewrite the text of the first run (before the bookmark) ---
$targetRange = $d.Paragraphs($targetIdx).Range
$oldRun1 = "MVN project (//vacsi problem ako sa zda; tie JRI kniznice a" + $nbsp + "spol. nemaju moc mvn repo a tak"
$newRun1 = "po nakresleni plotu sa prepnut na tab (takze po "
$ok1 = $targetRange.Find.Execute($oldRun1, $false, $false, $false, $false, $false, $true, 1, $false, $newRun1, 2)
if (-not $ok1) { throw "Could not replace first run text" }

# --- Step 2: rewrite the text of the second run (after the bookmark) ---
# Re-fetch the (now shorter) paragraph range and restrict the Find to it so
# we don't touch any of the many other ")" characters in the document.
$d = $word.ActiveDocument
$targetRange = $d.Paragraphs($targetIdx).Range
$ok2 = $targetRange.Find.Execute(")", $false, $false, $false, $false, $false, $true, 1, $false, "Plot selected, alebo po Run, ...)", 2)
if (-not $ok2) { throw "Could not replace second run text" }

# --- Step 3: append a brand-new paragraph after it with the original text,
#     split across two runs, via a raw OOXML fragment so the run boundary is
#     preserved exactly (instead of being re-coalesced by plain typing). ---
$d = $word.ActiveDocument
$targetRange = $d.Paragraphs($targetIdx).Range
$insertAt = $d.Range($targetRange.End - 1, $targetRange.End - 1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>'
$run1 = '<w:r>' + $rPr + '<w:t>MVN project</w:t></w:r>'
$run2 = '<w:r>' + $rPr + '<w:t xml:space="preserve"> (//vacsi problem ako sa zda; tie JRI kniznice a' + $nbsp + 'spol. nemaju moc mvn repo a tak)</w:t></w:r>'
$pPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/>' + $rPr + '</w:pPr>'
$newParaXml = '<w:p ' + $wNs + '>' + $pPr + $run1 + $run2 + '</w:p>'
$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertAt.InsertXML($frag)
